$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Range("F2").Value = 30.05449223518372
$ws.Range("F3").Value = 29.76475262641907
$ws.Range("F4").Value = 29.71687602996826
$ws.Range("F5").Value = 29.84806704521179
$ws.Range("F6").Value = 29.9085111618042
$ws.Range("F7").Value = 29.66334533691406
$ws.Range("F8").Value = 29.74172282218933
$ws.Range("F9").Value = 29.65336775779724
$ws.Range("F10").Value = 29.65936231613159
$ws.Range("F11").Value = 29.97253847122192
$ws.Range("F12").Value = 29.7814199924469
$ws.Range("F13").Value = 29.84281253814697
$ws.Range("F14").Value = 29.71160292625427
$ws.Range("F15").Value = 29.49990725517273
$ws.Range("F16").Value = 29.8912980556488
$ws.Range("F17").Value = 29.78654170036316
$ws.Range("F18").Value = 30.1444206237793
$ws.Range("F19").Value = 30.06649613380432
$ws.Range("F20").Value = 29.84225177764893
$ws.Range("F21").Value = 30.07104659080505

$ws = $wb.Worksheets.Item("run_2")
$ws.Range("F2").Value = 30.10622096061707
$ws.Range("F3").Value = 29.69258618354797
$ws.Range("F4").Value = 29.83807063102722
$ws.Range("F5").Value = 29.78295588493347
$ws.Range("F6").Value = 29.74206638336182
$ws.Range("F7").Value = 29.75846648216248
$ws.Range("F8").Value = 29.66606664657593
$ws.Range("F9").Value = 29.88372159004211
$ws.Range("F10").Value = 29.75860023498535
$ws.Range("F11").Value = 30.08910799026489
$ws.Range("F12").Value = 29.75185489654541
$ws.Range("F13").Value = 29.73409652709961
$ws.Range("F14").Value = 29.85681247711182
$ws.Range("F15").Value = 29.80641913414001
$ws.Range("F16").Value = 29.91411375999451
$ws.Range("F17").Value = 29.73997092247009
$ws.Range("F18").Value = 29.82440233230591
$ws.Range("F19").Value = 29.81028389930725
$ws.Range("F20").Value = 29.9933865070343
$ws.Range("F21").Value = 30.06592917442322

$ws = $wb.Worksheets.Item("run_3")
$ws.Range("F2").Value = 30.1581449508667
$ws.Range("F3").Value = 29.89680957794189
$ws.Range("F4").Value = 30.03415036201477
$ws.Range("F5").Value = 30.16022562980652
$ws.Range("F6").Value = 29.90883898735046
$ws.Range("F7").Value = 29.80102467536926
$ws.Range("F8").Value = 29.73245596885681
$ws.Range("F9").Value = 29.81217193603516
$ws.Range("F10").Value = 29.87986874580384
$ws.Range("F11").Value = 30.12013053894043
$ws.Range("F12").Value = 29.77164149284363
$ws.Range("F13").Value = 29.75386810302734
$ws.Range("F14").Value = 29.95186042785645
$ws.Range("F15").Value = 29.73042511940002
$ws.Range("F16").Value = 29.83797287940979
$ws.Range("F17").Value = 29.67706799507141
$ws.Range("F18").Value = 29.78482794761658
$ws.Range("F19").Value = 29.81416773796081
$ws.Range("F20").Value = 29.81194734573364
$ws.Range("F21").Value = 30.17433309555054

$ws = $wb.Worksheets.Item("run_4")
$ws.Range("F2").Value = 30.16321134567261
$ws.Range("F3").Value = 29.87343120574951
$ws.Range("F4").Value = 29.69489192962646
$ws.Range("F5").Value = 29.66487193107605
$ws.Range("F6").Value = 29.8112576007843
$ws.Range("F7").Value = 29.70791459083557
$ws.Range("F8").Value = 29.74257850646973
$ws.Range("F9").Value = 29.68389797210693
$ws.Range("F10").Value = 29.81452178955078
$ws.Range("F11").Value = 29.91542148590088
$ws.Range("F12").Value = 29.95841312408448
$ws.Range("F13").Value = 29.88687086105347
$ws.Range("F14").Value = 29.92565393447876
$ws.Range("F15").Value = 29.71771788597107
$ws.Range("F16").Value = 29.90666437149048
$ws.Range("F17").Value = 30.24825930595398
$ws.Range("F18").Value = 30.6692214012146
$ws.Range("F19").Value = 30.40072727203369
$ws.Range("F20").Value = 30.29193639755249
$ws.Range("F21").Value = 30.73220729827881

$ws = $wb.Worksheets.Item("run_5")
$ws.Range("F2").Value = 30.66563391685486
$ws.Range("F3").Value = 30.25231552124023
$ws.Range("F4").Value = 30.38260698318481
$ws.Range("F5").Value = 30.4077033996582
$ws.Range("F6").Value = 30.54141926765442
$ws.Range("F7").Value = 30.42685651779175
$ws.Range("F8").Value = 30.5037100315094
$ws.Range("F9").Value = 30.2572112083435
$ws.Range("F10").Value = 30.34380197525024
$ws.Range("F11").Value = 30.76593852043152
$ws.Range("F12").Value = 29.51127576828003
$ws.Range("F13").Value = 29.46178722381592
$ws.Range("F14").Value = 29.51126170158386
$ws.Range("F15").Value = 30.14964056015014
$ws.Range("F16").Value = 31.37946653366089
$ws.Range("F17").Value = 31.21654462814331
$ws.Range("F18").Value = 30.18591141700745
$ws.Range("F19").Value = 30.03041005134583
$ws.Range("F20").Value = 30.12016773223877
$ws.Range("F21").Value = 30.54921078681945
